$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 404.63635
$ws.Range("I4").Value = 221.66667
$ws.Range("J4").Value = 624.2
$ws.Range("K4").Value = 221.66667
$ws.Range("L4").Value = 624.2
$ws.Range("M4").Value = -107.66667
$ws.Range("N4").Value = -852.2
$ws.Range("H9").Value = 103.76923
$ws.Range("I9").Value = 104.916664
$ws.Range("J9").Value = 90
$ws.Range("K9").Value = 104.916664
$ws.Range("L9").Value = 90
$ws.Range("M9").Value = 64.083336
$ws.Range("N9").Value = -428
$ws.Range("H130").Value = 50460
$ws.Range("J130").Value = 50460
$ws.Range("L130").Value = 50460
$ws.Range("N130").Value = -60500
$ws.Range("H138").Value = 3846.1724
$ws.Range("I138").Value = 3033.9
$ws.Range("J138").Value = 4273.684
$ws.Range("K138").Value = 9101.700000000001
$ws.Range("L138").Value = 12821.052
$ws.Range("M138").Value = -3961.700000000001
$ws.Range("N138").Value = -23101.052

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H3").Value = 21741918
$ws.Range("I3").Value = 142858740
$ws.Range("K3").Value = 142858740
$ws.Range("M3").Value = -142858625
$ws.Range("H4").Value = 137
$ws.Range("I4").Value = 137
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 137
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = -21
$ws.Range("N4").ClearContents()
$ws.Range("H5").Value = 176.4
$ws.Range("J5").Value = 301
$ws.Range("L5").Value = 301
$ws.Range("N5").Value = -525
$ws.Range("H23").Value = 17500
$ws.Range("J23").Value = 17500
$ws.Range("L23").Value = 17500
$ws.Range("N23").Value = -18018
$ws.Range("H45").Value = 68616.734
$ws.Range("I45").Value = 84719.75
$ws.Range("K45").Value = 84719.75
$ws.Range("M45").Value = -84342.75
$ws.Range("H107").Value = 17450
$ws.Range("J107").Value = 17450
$ws.Range("L107").Value = 17450
$ws.Range("N107").Value = -25130
$ws.Range("H117").Value = 31000
$ws.Range("J117").Value = 31000
$ws.Range("L117").Value = 31000
$ws.Range("N117").Value = -40178
$ws.Range("H122").Value = 25002210
$ws.Range("I122").Value = 45456316
$ws.Range("J122").Value = 2744.3333
$ws.Range("K122").Value = 136368948
$ws.Range("L122").Value = 8232.999899999999
$ws.Range("M122").Value = -136366498
$ws.Range("N122").Value = -13132.9999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 176.4
$ws.Range("J4").Value = 301
$ws.Range("L4").Value = 301
$ws.Range("N4").Value = -531
$ws.Range("H12").Value = 2125
$ws.Range("I12").Value = 1000
$ws.Range("J12").Value = 10000
$ws.Range("K12").Value = 1000
$ws.Range("L12").Value = 10000
$ws.Range("M12").Value = -832
$ws.Range("N12").Value = -10336
$ws.Range("H23").Value = 10000
$ws.Range("I23").Value = 0
$ws.Range("J23").Value = 10000
$ws.Range("K23").Value = 0
$ws.Range("L23").Value = 10000
$ws.Range("M23").ClearContents()
$ws.Range("N23").Value = -10566
$ws.Range("H105").Value = 2777.5
$ws.Range("I105").Value = 2536.6667
$ws.Range("J105").Value = 3500
$ws.Range("K105").Value = 2536.6667
$ws.Range("L105").Value = 3500
$ws.Range("M105").Value = -789.6667000000002
$ws.Range("N105").Value = -6994
$ws.Range("H134").Value = 29424.514
$ws.Range("I134").Value = 33416.53
$ws.Range("K134").Value = 100249.59
$ws.Range("M134").Value = -97714.59
$ws.Range("H135").Value = 57569.855
$ws.Range("J135").Value = 57569.855
$ws.Range("L135").Value = 57569.855
$ws.Range("N135").Value = -67709.85500000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H6").Value = 0
$ws.Range("I6").Value = 0
$ws.Range("J6").Value = 0
$ws.Range("K6").Value = 0
$ws.Range("L6").Value = 0
$ws.Range("M6").ClearContents()
$ws.Range("N6").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 938
$ws.Range("I3").Value = 938
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 2814
$ws.Range("L3").Value = 0
$ws.Range("M3").Value = -2702
$ws.Range("N3").ClearContents()
$ws.Range("H31").Value = 0
$ws.Range("J31").Value = 0
$ws.Range("L31").Value = 0
$ws.Range("N31").ClearContents()
$ws.Range("H35").Value = 2000
$ws.Range("J35").Value = 2000
$ws.Range("L35").Value = 6000
$ws.Range("N35").Value = -6576
$ws.Range("H60").Value = 30304458
$ws.Range("J60").Value = 3125
$ws.Range("L60").Value = 9375
$ws.Range("N60").Value = -9877
$ws.Range("H100").Value = 2240
$ws.Range("J100").Value = 2500
$ws.Range("L100").Value = 7500
$ws.Range("N100").Value = -9122
$ws.Range("H114").Value = 1510.1666
$ws.Range("J114").Value = 2327.9333
$ws.Range("L114").Value = 6983.7999
$ws.Range("N114").Value = -13491.7999
$ws.Range("H117").Value = 2900
$ws.Range("I117").Value = 700
$ws.Range("J117").Value = 4000
$ws.Range("K117").Value = 2100
$ws.Range("L117").Value = 12000
$ws.Range("M117").Value = 1342
$ws.Range("N117").Value = -18884
$ws.Range("H125").Value = 2660
$ws.Range("I125").Value = 3533.3333
$ws.Range("J125").Value = 2285.7144
$ws.Range("K125").Value = 10599.9999
$ws.Range("L125").Value = 6857.1432
$ws.Range("M125").Value = -5679.999899999999
$ws.Range("N125").Value = -16697.1432
$ws.Range("H131").Value = 1221201.9
$ws.Range("J131").Value = 1300352.2
$ws.Range("L131").Value = 3901056.6
$ws.Range("N131").Value = -3911136.6

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H45").Value = 0
$ws.Range("J45").Value = 0
$ws.Range("L45").Value = 0
$ws.Range("N45").ClearContents()
$ws.Range("H51").Value = 43632.168
$ws.Range("J51").Value = 43871.453
$ws.Range("L51").Value = 43871.453
$ws.Range("N51").Value = -44889.453
$ws.Range("H102").Value = 1555.4286
$ws.Range("I102").Value = 1555.4286
$ws.Range("K102").Value = 1555.4286
$ws.Range("M102").Value = 66.57140000000004
$ws.Range("H126").Value = 1648.4
$ws.Range("I126").Value = 1447.3334
$ws.Range("K126").Value = 4342.0002
$ws.Range("M126").Value = -1872.0002

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 2607
$ws.Range("I40").Value = 2041.5
$ws.Range("K40").Value = 2041.5
$ws.Range("M40").Value = -1905.5
$ws.Range("H58").Value = 2000
$ws.Range("I58").Value = 2000
$ws.Range("K58").Value = 2000
$ws.Range("M58").Value = -1740
$ws.Range("H61").Value = 2755.7778
$ws.Range("I61").Value = 1543.4286
$ws.Range("J61").Value = 6999
$ws.Range("K61").Value = 1543.4286
$ws.Range("L61").Value = 6999
$ws.Range("M61").Value = -1341.4286
$ws.Range("N61").Value = -7403
$ws.Range("H113").Value = 2755.7778
$ws.Range("I113").Value = 1543.4286
$ws.Range("J113").Value = 6999
$ws.Range("K113").Value = 1543.4286
$ws.Range("L113").Value = 6999
$ws.Range("M113").Value = 626.5714
$ws.Range("N113").Value = -11339
$ws.Range("H122").Value = 21175.875
$ws.Range("I122").Value = 32201.4
$ws.Range("K122").Value = 96604.20000000001
$ws.Range("M122").Value = -94154.20000000001
$ws.Range("H132").Value = 15546.556
$ws.Range("I132").Value = 20336.666
$ws.Range("J132").Value = 5966.3335
$ws.Range("K132").Value = 61009.99800000001
$ws.Range("L132").Value = 17899.0005
$ws.Range("M132").Value = -58479.99800000001
$ws.Range("N132").Value = -22959.0005

$wb.Save()